# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Modelo de Temuco" (Mango)
# as row 195, pushing the existing rows 195-198 down to 196-199.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 195:198 down by inserting a new blank row at 195.
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new weekly record.
$ws.Cells.Item(195, 1).Value  = 10
$ws.Cells.Item(195, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(195, 3).Value  = "La Araucanía"
$ws.Cells.Item(195, 4).Value  = 44448
$ws.Cells.Item(195, 5).Value  = 9
$ws.Cells.Item(195, 6).Value  = "Fruta"
$ws.Cells.Item(195, 7).Value  = 100108
$ws.Cells.Item(195, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(195, 9).Value  = 100108002
$ws.Cells.Item(195, 10).Value = "Mango"
$ws.Cells.Item(195, 11).Value = "Sin especificar"
$ws.Cells.Item(195, 12).Value = "Primera"
$ws.Cells.Item(195, 13).Value = 800
$ws.Cells.Item(195, 14).Value = 9000
$ws.Cells.Item(195, 15).Value = 9000
$ws.Cells.Item(195, 16).Value = 9000
$ws.Cells.Item(195, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(195, 18).Value = "Brasil"
$ws.Cells.Item(195, 19).Value = 2250
$ws.Cells.Item(195, 20).Value = 4
